$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 117 (shifts existing rows 117-125 down to 118-126)
$ws.Rows.Item(117).Insert()

# Populate the new row 117 with the new weekly price record
$ws.Cells.Item(117, 1).Value = 3
$ws.Cells.Item(117, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(117, 3).Value = "Coquimbo"
$ws.Cells.Item(117, 4).Value = 44585
$ws.Cells.Item(117, 5).Value = 5
$ws.Cells.Item(117, 6).Value = 100112052
$ws.Cells.Item(117, 7).Value = "Albahaca"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 70
$ws.Cells.Item(117, 11).Value = 4000
$ws.Cells.Item(117, 12).Value = 4500
$ws.Cells.Item(117, 13).Value = 4286
$ws.Cells.Item(117, 14).Value = "`$/docena de matas"
$ws.Cells.Item(117, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(117, 16).Value = 714
$ws.Cells.Item(117, 17).Value = 6
$ws.Cells.Item(117, 18).Value = "Hortaliza"
